# Scheduled market-data refresh: updates the computed price/profit columns
# (H:N) on each job sheet. Values come from an external price feed, so
# cells are written as plain numbers (no formulas). A few cells that no
# longer have a computed value are cleared entirely rather than zeroed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 459.85715
$ws.Range("I28").Value = 395.81818
$ws.Range("K28").Value = 395.81818
$ws.Range("M28").Value = 89.18182000000002

$ws.Range("H51").Value = 83354170
$ws.Range("I51").Value = 250012500
$ws.Range("J51").Value = 24999.25
$ws.Range("K51").Value = 250012500
$ws.Range("L51").Value = 24999.25
$ws.Range("M51").Value = -250012016
$ws.Range("N51").Value = -25967.25

$ws.Range("H111").Value = 3235.2666
$ws.Range("I111").Value = 4854.875
$ws.Range("J111").Value = 1384.2858
$ws.Range("K111").Value = 14564.625
$ws.Range("L111").Value = 4152.857400000001
$ws.Range("M111").Value = -11497.625
$ws.Range("N111").Value = -10286.8574

$ws.Range("H116").Value = 9998.799999999999
$ws.Range("J116").Value = 14666.667
$ws.Range("L116").Value = 14666.667
$ws.Range("N116").Value = -21550.667

$ws.Range("H125").Value = 937.1111
$ws.Range("I125").Value = 929.375
$ws.Range("K125").Value = 8364.375
$ws.Range("M125").Value = -5904.375

$ws.Range("H141").Value = 3836.6667
$ws.Range("I141").Value = 3566.25
$ws.Range("J141").Value = 6000
$ws.Range("K141").Value = 10698.75
$ws.Range("L141").Value = 18000
$ws.Range("M141").Value = -5518.75
$ws.Range("N141").Value = -28360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H74").Value = 2030
$ws.Range("I74").Value = 1241
$ws.Range("J74").Value = 3044.4285
$ws.Range("K74").Value = 1241
$ws.Range("L74").Value = 3044.4285
$ws.Range("M74").Value = -367
$ws.Range("N74").Value = -4792.4285

$ws.Range("H77").Value = 2030
$ws.Range("I77").Value = 1241
$ws.Range("J77").Value = 3044.4285
$ws.Range("K77").Value = 6205
$ws.Range("L77").Value = 15222.1425
$ws.Range("M77").Value = -1837
$ws.Range("N77").Value = -23958.1425

$ws.Range("H132").Value = 2585.0952
$ws.Range("I132").Value = 2878.1
$ws.Range("J132").Value = 2318.7273
$ws.Range("K132").Value = 8634.299999999999
$ws.Range("L132").Value = 6956.1819
$ws.Range("M132").Value = -6104.299999999999
$ws.Range("N132").Value = -12016.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 26197
$ws.Range("J81").Value = 27945
$ws.Range("L81").Value = 27945
$ws.Range("N81").Value = -30067

$ws.Range("H84").Value = 26197
$ws.Range("J84").Value = 27945
$ws.Range("L84").Value = 83835
$ws.Range("N84").Value = -94443

$ws.Range("H86").Value = 5234.6665
$ws.Range("I86").Value = 5131.8
$ws.Range("K86").Value = 5131.8
$ws.Range("M86").Value = -4008.8

$ws.Range("H89").Value = 5234.6665
$ws.Range("I89").Value = 5131.8
$ws.Range("K89").Value = 25659
$ws.Range("M89").Value = -20043

$ws.Range("H105").Value = 11819741
$ws.Range("I105").Value = 1001279.7
$ws.Range("J105").Value = 20835124
$ws.Range("K105").Value = 1001279.7
$ws.Range("L105").Value = 20835124
$ws.Range("M105").Value = -999532.7
$ws.Range("N105").Value = -20838618

$ws.Range("H107").Value = 2137978.8
$ws.Range("J107").Value = 1499.1111
$ws.Range("L107").Value = 1499.1111
$ws.Range("N107").Value = -5339.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1938.6364
$ws.Range("I58").Value = 1386.5
$ws.Range("J58").Value = 2398.75
$ws.Range("K58").Value = 1386.5
$ws.Range("L58").Value = 2398.75
$ws.Range("M58").Value = -1183.5
$ws.Range("N58").Value = -2804.75

$ws.Range("H69").Value = 36999.4
$ws.Range("I69").Value = 34998.5
$ws.Range("K69").Value = 34998.5
$ws.Range("M69").Value = -34249.5

$ws.Range("H72").Value = 36999.4
$ws.Range("I72").Value = 34998.5
$ws.Range("K72").Value = 104995.5
$ws.Range("M72").Value = -101251.5

$ws.Range("H132").Value = 5037
$ws.Range("I132").Value = 4598.222
$ws.Range("J132").Value = 6024.25
$ws.Range("K132").Value = 13794.666
$ws.Range("L132").Value = 18072.75
$ws.Range("M132").Value = -11264.666
$ws.Range("N132").Value = -23132.75

$ws.Range("H134").Value = 5646.263
$ws.Range("I134").Value = 5861.1875
$ws.Range("K134").Value = 17583.5625
$ws.Range("M134").Value = -15048.5625

$ws.Range("H136").Value = 1938.6364
$ws.Range("I136").Value = 1386.5
$ws.Range("J136").Value = 2398.75
$ws.Range("K136").Value = 4159.5
$ws.Range("L136").Value = 7196.25
$ws.Range("M136").Value = -1609.5
$ws.Range("N136").Value = -12296.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 475
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws.Range("H33").Value = 280
$ws.Range("I33").Value = 280
$ws.Range("J33").Value = 280
$ws.Range("K33").Value = 1680
$ws.Range("L33").Value = 1680
$ws.Range("M33").Value = -1397
$ws.Range("N33").Value = -2246

$ws.Range("H68").Value = 6671630.5
$ws.Range("J68").Value = 11118666
$ws.Range("L68").Value = 33355998
$ws.Range("N68").Value = -33357620

$ws.Range("H71").Value = 6671630.5
$ws.Range("J71").Value = 11118666
$ws.Range("L71").Value = 100067994
$ws.Range("N71").Value = -100076106

$ws.Range("H107").Value = 491.0909

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1148.9524
$ws.Range("I97").Value = 1215.7778
$ws.Range("J97").Value = 748
$ws.Range("K97").Value = 1215.7778
$ws.Range("L97").Value = 748
$ws.Range("M97").Value = -719.7778000000001
$ws.Range("N97").Value = -1740

$ws.Range("H122").Value = 3245.4348
$ws.Range("I122").Value = 2156.1667
$ws.Range("K122").Value = 6468.500100000001
$ws.Range("M122").Value = -4018.500100000001

$ws.Range("H126").Value = 8769.951999999999
$ws.Range("I126").Value = 2138.4
$ws.Range("J126").Value = 10842.3125
$ws.Range("K126").Value = 6415.200000000001
$ws.Range("L126").Value = 32526.9375
$ws.Range("M126").Value = -3945.200000000001
$ws.Range("N126").Value = -37466.9375

$ws.Range("H132").Value = 1798.8823
$ws.Range("I132").Value = 1468.3846
$ws.Range("K132").Value = 4405.1538
$ws.Range("M132").Value = -1875.1538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2537.375
$ws.Range("I7").Value = 1576
$ws.Range("J7").Value = 3498.75
$ws.Range("K7").Value = 1576
$ws.Range("L7").Value = 3498.75
$ws.Range("M7").Value = -1464
$ws.Range("N7").Value = -3722.75

$ws.Range("H68").Value = 1774.5714
$ws.Range("I68").Value = 1670.3334
$ws.Range("K68").Value = 1670.3334
$ws.Range("M68").Value = -921.3334

$ws.Range("H71").Value = 1774.5714
$ws.Range("I71").Value = 1670.3334
$ws.Range("K71").Value = 8351.666999999999
$ws.Range("M71").Value = -4607.666999999999

$ws.Range("H126").Value = 2537.375
$ws.Range("I126").Value = 1576
$ws.Range("J126").Value = 3498.75
$ws.Range("K126").Value = 4728
$ws.Range("L126").Value = 10496.25
$ws.Range("M126").Value = -2258
$ws.Range("N126").Value = -15436.25

$ws.Range("H132").Value = 4583.2856
$ws.Range("I132").Value = 5134.5454
$ws.Range("K132").Value = 15403.6362
$ws.Range("M132").Value = -12873.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4152.3687
$ws.Range("I132").Value = 4502.2856
$ws.Range("K132").Value = 13506.8568
$ws.Range("M132").Value = -10976.8568
